# Insert a new weekly record at row 33 (Región del Maule, serial date 44526
# = 2021-11-26), pushing the existing rows 33-111 down to 34-112.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(33).Insert()

$ws.Range("A33").Value2 = 5
$ws.Range("B33").Value2 = "Macroferia Regional de Talca"
$ws.Range("C33").Value2 = "Maule"
$ws.Range("D33").Value2 = 44526
$ws.Range("E33").Value2 = 7
$ws.Range("F33").Value2 = 100112031
$ws.Range("G33").Value2 = "Poroto verde"
$ws.Range("H33").Value2 = "Sin especificar"
$ws.Range("I33").Value2 = "Primera"
$ws.Range("J33").Value2 = 150
$ws.Range("K33").Value2 = 30000
$ws.Range("L33").Value2 = 30000
$ws.Range("M33").Value2 = 30000
$ws.Range("N33").Value2 = "`$/saco 25 kilos"
$ws.Range("O33").Value2 = "Región del Maule"
$ws.Range("P33").Value2 = 1200
$ws.Range("Q33").Value2 = 25
$ws.Range("R33").Value2 = "Hortaliza"
